$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Weight full Falcon tube" record (row 22), which is a superfluous
# weight determination. All following rows shift up by one.
$ws.Rows.Item(22).Delete()
